$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temperature ramp data Run 5 - updated raw ramp-start readings.
# E4 and E16 are entered (non-formula) anchor values; E7 is overwritten
# with a literal reading (no longer computed); E19 becomes formula-driven
# again so it continues the E5:E24 running-total chain.
$ws.Range("E4").Value = 0.56597222222222221
$ws.Range("E7").Value = 0.59930555555555554
$ws.Range("E16").Value = 0.70000000000000007
$ws.Range("E19").Formula = "=E18+D19"

# Selection left on E17 after entering the new readings.
[void]$ws.Range("E17").Select()
